# Updated cryptos list on Sat Jun 10 15:09:07 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns of the cryptos sheet
# with the latest coinranking.com scrape. Rows 24/25 (Chainlink/Cosmos) and
# 42/43 (Quant/VeChain) also swap rank position, so Coin (B) and Link (C)
# are updated for those rows as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '25.674.25'
$ws.Range("E2").Value = '  -3.54%  '
# Row 3
$ws.Range("D3").Value = '1.743.73'
$ws.Range("E3").Value = '  -5.65%  '
# Row 4
$ws.Range("E4").Value = '  +0.07%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.24'
$ws.Range("E5").Value = '  -10.42%  '
# Row 6
$ws.Range("E6").Value = '  -0.06%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4923'
$ws.Range("E7").Value = '  -7.89%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.61'
$ws.Range("E8").Value = '  -7.81%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2547'
$ws.Range("E9").Value = '  -19.54%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06009'
$ws.Range("E10").Value = '  -13.59%  '
# Row 11
$ws.Range("D11").Value = '1.743.82'
$ws.Range("E11").Value = '  -5.81%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.06843'
$ws.Range("E12").Value = '  -12.70%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.74'
$ws.Range("E13").Value = '  -21.86%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.444'
$ws.Range("E14").Value = '  -11.93%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '76.71'
$ws.Range("E15").Value = '  -14.38%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.5650'
$ws.Range("E16").Value = '  -26.66%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.002'
$ws.Range("E17").Value = '  +0.07%  '
# Row 18
$ws.Range("E18").Value = '  -0.02%  '
# Row 19
$ws.Range("D19").Value = '25.724.64'
$ws.Range("E19").Value = '  -3.47%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.25'
$ws.Range("E20").Value = '  -20.36%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.000006570'
$ws.Range("E21").Value = '  -17.67%  '
# Row 22
$ws.Range("D22").Value = '1.967.94'
$ws.Range("E22").Value = '  -6.00%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.005'
$ws.Range("E23").Value = '  -13.77%  '
# Row 24
$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.879'
$ws.Range("E24").Value = '  -15.86%  '
# Row 25
$ws.Range("B25").Value = 'Chainlink'
$ws.Range("C25").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.015'
$ws.Range("E25").Value = '  -16.70%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '136.97'
$ws.Range("E26").Value = '  -3.56%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.488'
$ws.Range("E27").Value = '  -12.02%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.818'
$ws.Range("E28").Value = '  -17.90%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '14.68'
$ws.Range("E29").Value = '  -14.29%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '101.89'
$ws.Range("E30").Value = '  -8.77%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.757'
$ws.Range("E31").Value = '  -12.79%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.07990'
$ws.Range("E32").Value = '  -8.98%  '
# Row 33
$ws.Range("E33").Value = '  -17.96%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04390'
$ws.Range("E34").Value = '  -9.57%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9998'
$ws.Range("E35").Value = '  -0.13%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.602'
$ws.Range("E36").Value = '  -9.84%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9774'
$ws.Range("E37").Value = '  -14.18%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6014'
$ws.Range("E38").Value = '  -18.35%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.668'
$ws.Range("E39").Value = '  -14.08%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.995'
$ws.Range("E40").Value = '  -15.41%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.000'
$ws.Range("E41").Value = '  -0.10%  '
# Row 42
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01501'
$ws.Range("E42").Value = '  -13.60%  '
# Row 43
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '102.44'
$ws.Range("E43").Value = '  -6.01%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7528'
$ws.Range("E44").Value = '  -17.16%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.158'
$ws.Range("E45").Value = '  -12.73%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.3727'
$ws.Range("E46").Value = '  -22.75%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05238'
$ws.Range("E47").Value = '  -9.97%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1062'
$ws.Range("E48").Value = '  -15.06%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.00'
$ws.Range("E49").Value = '  -14.40%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.856'
$ws.Range("E50").Value = '  -24.00%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '52.19'
$ws.Range("E51").Value = '  -13.56%  '
